$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 8677.714
$ws.Range("I32").Value = 25500.5
$ws.Range("K32").Value = 25500.5
$ws.Range("M32").Value = -25174.5
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H69").Value = 7688.3335
$ws.Range("I69").Value = 5550
$ws.Range("K69").Value = 16650
$ws.Range("M69").Value = -15776
$ws.Range("H72").Value = 7688.3335
$ws.Range("I72").Value = 5550
$ws.Range("K72").Value = 49950
$ws.Range("M72").Value = -45582
$ws.Range("H98").Value = 33001.2
$ws.Range("J98").Value = 3006
$ws.Range("L98").Value = 3006
$ws.Range("N98").Value = -6002
$ws.Range("H113").Value = 2668.2
$ws.Range("J113").Value = 3141.2
$ws.Range("L113").Value = 3141.2
$ws.Range("N113").Value = -9649.200000000001
$ws.Range("H116").Value = 7912.6
$ws.Range("I116").Value = 6624
$ws.Range("J116").Value = 9201.200000000001
$ws.Range("K116").Value = 6624
$ws.Range("L116").Value = 9201.200000000001
$ws.Range("M116").Value = -3182
$ws.Range("N116").Value = -16085.2
$ws.Range("H121").Value = 1533.3334
$ws.Range("J121").Value = 1533.3334
$ws.Range("L121").Value = 4600.0002
$ws.Range("N121").Value = -8094.0002
$ws.Range("H122").Value = 33001.2
$ws.Range("J122").Value = 3006
$ws.Range("L122").Value = 9018
$ws.Range("N122").Value = -13918
$ws.Range("H132").Value = 1573.9286
$ws.Range("I132").Value = 1114.84
$ws.Range("K132").Value = 3344.52
$ws.Range("M132").Value = -814.5199999999995
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 36000.676
$ws.Range("I32").Value = 39065.742
$ws.Range("K32").Value = 39065.742
$ws.Range("M32").Value = -38778.742
$ws.Range("H45").Value = 1462.8948
$ws.Range("I45").Value = 1318.5625
$ws.Range("K45").Value = 1318.5625
$ws.Range("M45").Value = -941.5625
$ws.Range("H61").Value = 12829584
$ws.Range("I61").Value = 20839424
$ws.Range("J61").Value = 13841.5
$ws.Range("K61").Value = 20839424
$ws.Range("L61").Value = 13841.5
$ws.Range("M61").Value = -20839212
$ws.Range("N61").Value = -14265.5
$ws.Range("H69").Value = 338665.6
$ws.Range("J69").Value = 338665.6
$ws.Range("L69").Value = 338665.6
$ws.Range("N69").Value = -340163.6
$ws.Range("H72").Value = 338665.6
$ws.Range("J72").Value = 338665.6
$ws.Range("L72").Value = 1015996.8
$ws.Range("N72").Value = -1023484.8
$ws.Range("H74").Value = 3193.842
$ws.Range("I74").Value = 1315.2222
$ws.Range("K74").Value = 1315.2222
$ws.Range("M74").Value = -441.2221999999999
$ws.Range("H77").Value = 3193.842
$ws.Range("I77").Value = 1315.2222
$ws.Range("K77").Value = 6576.111
$ws.Range("M77").Value = -2208.111
$ws.Range("H132").Value = 3778102
$ws.Range("I132").Value = 5558437.5
$ws.Range("J132").Value = 7979.5884
$ws.Range("K132").Value = 16675312.5
$ws.Range("L132").Value = 23938.7652
$ws.Range("M132").Value = -16672782.5
$ws.Range("N132").Value = -28998.7652
$ws.Range("H136").Value = 12829584
$ws.Range("I136").Value = 20839424
$ws.Range("J136").Value = 13841.5
$ws.Range("K136").Value = 62518272
$ws.Range("L136").Value = 41524.5
$ws.Range("M136").Value = -62515722
$ws.Range("N136").Value = -46624.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 83356060
$ws.Range("I105").Value = 125032424
$ws.Range("K105").Value = 125032424
$ws.Range("M105").Value = -125030677
$ws.Range("H107").Value = 1375.0667
$ws.Range("I107").Value = 1030.2174
$ws.Range("K107").Value = 1030.2174
$ws.Range("M107").Value = 889.7826
$ws.Range("H118").Value = 57499.668
$ws.Range("J118").Value = 57499.668
$ws.Range("L118").Value = 57499.668
$ws.Range("N118").Value = -60813.668
$ws.Range("H131").Value = 54625
$ws.Range("J131").Value = 54625
$ws.Range("L131").Value = 54625
$ws.Range("N131").Value = -64705
$ws.Range("H134").Value = 6187.706
$ws.Range("I134").Value = 4817.45
$ws.Range("K134").Value = 14452.35
$ws.Range("M134").Value = -11917.35
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 128.58824
$ws.Range("I7").Value = 145.3077
$ws.Range("J7").Value = 74.25
$ws.Range("K7").Value = 145.3077
$ws.Range("L7").Value = 74.25
$ws.Range("M7").Value = -32.30770000000001
$ws.Range("N7").Value = -300.25
$ws.Range("H69").Value = 13499
$ws.Range("I69").Value = 7000
$ws.Range("J69").Value = 19998
$ws.Range("K69").Value = 7000
$ws.Range("L69").Value = 19998
$ws.Range("N69").Value = -21496
$ws.Range("M69").Value = -6251
$ws.Range("H72").Value = 13499
$ws.Range("I72").Value = 7000
$ws.Range("J72").Value = 19998
$ws.Range("K72").Value = 21000
$ws.Range("L72").Value = 59994
$ws.Range("N72").Value = -67482
$ws.Range("M72").Value = -17256
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H99").Value = 8999.333000000001
$ws.Range("I99").Value = 7998
$ws.Range("K99").Value = 7998
$ws.Range("M99").Value = -6500
$ws.Range("H122").Value = 104640.9
$ws.Range("J122").Value = 6989.8
$ws.Range("L122").Value = 20969.4
$ws.Range("N122").Value = -25869.4
$ws.Range("H126").Value = 8999.333000000001
$ws.Range("I126").Value = 7998
$ws.Range("K126").Value = 23994
$ws.Range("M126").Value = -21524
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 4014.8
$ws.Range("I63").Value = 4515
$ws.Range("J63").Value = 2014
$ws.Range("K63").Value = 13545
$ws.Range("L63").Value = 6042
$ws.Range("M63").Value = -12796
$ws.Range("N63").Value = -7540
$ws.Range("H66").Value = 4014.8
$ws.Range("I66").Value = 4515
$ws.Range("J66").Value = 2014
$ws.Range("K66").Value = 40635
$ws.Range("L66").Value = 18126
$ws.Range("M66").Value = -36891
$ws.Range("N66").Value = -25614
$ws.Range("H70").Value = 3423.1667
$ws.Range("I70").Value = 3834.75
$ws.Range("K70").Value = 11504.25
$ws.Range("M70").Value = -11189.25
$ws.Range("H73").Value = 3423.1667
$ws.Range("I73").Value = 3834.75
$ws.Range("K73").Value = 11504.25
$ws.Range("M73").Value = -10412.25
$ws.Range("H101").Value = 6162.5
$ws.Range("J101").Value = 9999
$ws.Range("L101").Value = 29997
$ws.Range("N101").Value = -34865
$ws.Range("H138").Value = 911390.8
$ws.Range("I138").Value = 1525
$ws.Range("K138").Value = 4575
$ws.Range("M138").Value = 565
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3594.7917
$ws.Range("I22").Value = 1862
$ws.Range("K22").Value = 1862
$ws.Range("M22").Value = -1567
$ws.Range("H27").Value = 3594.7917
$ws.Range("I27").Value = 1862
$ws.Range("K27").Value = 1862
$ws.Range("M27").Value = -1755
$ws.Range("H40").Value = 14802.814
$ws.Range("I40").Value = 11746.714
$ws.Range("K40").Value = 11746.714
$ws.Range("M40").Value = -11610.714
$ws.Range("H46").Value = 8364.682000000001
$ws.Range("J46").Value = 8667.714
$ws.Range("L46").Value = 8667.714
$ws.Range("N46").Value = -9043.714
$ws.Range("H55").Value = 159.72223
$ws.Range("I55").Value = 150.85715
$ws.Range("K55").Value = 150.85715
$ws.Range("M55").Value = 22.14285000000001
$ws.Range("H82").Value = 1654.8182
$ws.Range("I82").Value = 1816.6666
$ws.Range("J82").Value = 1460.6
$ws.Range("K82").Value = 1816.6666
$ws.Range("L82").Value = 1460.6
$ws.Range("M82").Value = -1455.6666
$ws.Range("N82").Value = -2182.6
$ws.Range("H85").Value = 1654.8182
$ws.Range("I85").Value = 1816.6666
$ws.Range("J85").Value = 1460.6
$ws.Range("K85").Value = 1816.6666
$ws.Range("L85").Value = 1460.6
$ws.Range("M85").Value = -568.6666
$ws.Range("N85").Value = -3956.6
$ws.Range("H136").Value = 1195466.5
$ws.Range("I136").Value = 1485682.8
$ws.Range("K136").Value = 4457048.4
$ws.Range("M136").Value = -4454498.4
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4374.75
$ws.Range("I62").Value = 3999
$ws.Range("J62").Value = 4500
$ws.Range("K62").Value = 3999
$ws.Range("L62").Value = 4500
$ws.Range("M62").Value = -3375
$ws.Range("N62").Value = -5748
$ws.Range("H65").Value = 4374.75
$ws.Range("I65").Value = 3999
$ws.Range("J65").Value = 4500
$ws.Range("K65").Value = 19995
$ws.Range("L65").Value = 22500
$ws.Range("M65").Value = -16875
$ws.Range("N65").Value = -28740
$ws.Range("H74").Value = 14000
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 14000
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 14000
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -15872
$ws.Range("H77").Value = 14000
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 14000
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 42000
$ws.Range("M77").ClearContents()
